$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 58.833332
$ws.Range("J9").Value = 54.5
$ws.Range("L9").Value = 54.5
$ws.Range("N9").Value = -392.5

$ws.Range("H32").Value = 799.5
$ws.Range("J32").Value = 799.5
$ws.Range("L32").Value = 799.5
$ws.Range("N32").Value = -1451.5

$ws.Range("H53").Value = 459
$ws.Range("I53").Value = 551.4286
$ws.Range("K53").Value = 551.4286
$ws.Range("M53").Value = 85.57140000000004

$ws.Range("H55").Value = 832.2222
$ws.Range("J55").Value = 961.4545000000001
$ws.Range("L55").Value = 961.4545000000001
$ws.Range("N55").Value = -1389.4545

$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("M76").ClearContents()
$ws.Range("N76").ClearContents()

$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("M79").ClearContents()
$ws.Range("N79").ClearContents()

$ws.Range("H86").Value = 5249
$ws.Range("J86").Value = 9998
$ws.Range("L86").Value = 9998
$ws.Range("N86").Value = -12244

$ws.Range("H89").Value = 5249
$ws.Range("J89").Value = 9998
$ws.Range("L89").Value = 49990
$ws.Range("N89").Value = -61222

$ws.Range("H98").Value = 301.42856
$ws.Range("I98").Value = 301.42856
$ws.Range("K98").Value = 301.42856
$ws.Range("M98").Value = 1196.57144

$ws.Range("H116").Value = 2884.9285
$ws.Range("I116").Value = 2800
$ws.Range("J116").Value = 3097.25
$ws.Range("K116").Value = 2800
$ws.Range("L116").Value = 3097.25
$ws.Range("M116").Value = 642
$ws.Range("N116").Value = -9981.25

$ws.Range("H122").Value = 301.42856
$ws.Range("I122").Value = 301.42856
$ws.Range("K122").Value = 904.28568
$ws.Range("M122").Value = 1545.71432

$ws.Range("H135").Value = 356.42856
$ws.Range("I135").Value = 393.16666
$ws.Range("J135").Value = 136
$ws.Range("K135").Value = 3538.49994
$ws.Range("L135").Value = 1224
$ws.Range("M135").Value = -1003.49994
$ws.Range("N135").Value = -6294

$ws.Range("H137").Value = 2533.2273
$ws.Range("I137").Value = 1716.5
$ws.Range("K137").Value = 5149.5
$ws.Range("M137").Value = -2599.5

$ws.Range("H138").Value = 3078.6924
$ws.Range("I138").Value = 1638
$ws.Range("K138").Value = 4914
$ws.Range("M138").Value = 226

$ws.Range("H140").Value = 38926.668
$ws.Range("J140").Value = 38926.668
$ws.Range("L140").Value = 38926.668
$ws.Range("N140").Value = -49286.668

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 5299.6665
$ws.Range("I63").Value = 2200
$ws.Range("J63").Value = 11499
$ws.Range("K63").Value = 2200
$ws.Range("L63").Value = 11499
$ws.Range("M63").Value = -1514
$ws.Range("N63").Value = -12871

$ws.Range("H66").Value = 5299.6665
$ws.Range("I66").Value = 2200
$ws.Range("J66").Value = 11499
$ws.Range("K66").Value = 11000
$ws.Range("L66").Value = 57495
$ws.Range("M66").Value = -7568
$ws.Range("N66").Value = -64359

$ws.Range("H102").Value = 4474.5713
$ws.Range("I102").Value = 2079.75
$ws.Range("K102").Value = 2079.75
$ws.Range("M102").Value = -457.75

$ws.Range("H132").Value = 1888.1111
$ws.Range("I132").Value = 1956.1428
$ws.Range("K132").Value = 5868.428400000001
$ws.Range("M132").Value = -3338.428400000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 20000
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()

$ws.Range("H64").Value = 539.375
$ws.Range("I64").Value = 489.5
$ws.Range("J64").Value = 589.25
$ws.Range("K64").Value = 489.5
$ws.Range("L64").Value = 589.25
$ws.Range("M64").Value = -264.5
$ws.Range("N64").Value = -1039.25

$ws.Range("H67").Value = 539.375
$ws.Range("I67").Value = 489.5
$ws.Range("J67").Value = 589.25
$ws.Range("K67").Value = 489.5
$ws.Range("L67").Value = 589.25
$ws.Range("M67").Value = 290.5
$ws.Range("N67").Value = -2149.25

$ws.Range("H76").Value = 20000
$ws.Range("J76").Value = 20000
$ws.Range("L76").Value = 20000
$ws.Range("N76").Value = -20630

$ws.Range("H79").Value = 20000
$ws.Range("J79").Value = 20000
$ws.Range("L79").Value = 20000
$ws.Range("N79").Value = -22184

$ws.Range("H86").Value = 2754.7585
$ws.Range("I86").Value = 1669.7
$ws.Range("J86").Value = 5166
$ws.Range("K86").Value = 1669.7
$ws.Range("L86").Value = 5166
$ws.Range("M86").Value = -546.7
$ws.Range("N86").Value = -7412

$ws.Range("H89").Value = 2754.7585
$ws.Range("I89").Value = 1669.7
$ws.Range("J89").Value = 5166
$ws.Range("K89").Value = 8348.5
$ws.Range("L89").Value = 25830
$ws.Range("M89").Value = -2732.5
$ws.Range("N89").Value = -37062

$ws.Range("H107").Value = 4589.1763
$ws.Range("J107").Value = 7996.1665
$ws.Range("L107").Value = 7996.1665
$ws.Range("N107").Value = -11836.1665

$ws.Range("H134").Value = 1008.6667
$ws.Range("I134").Value = 1009.5
$ws.Range("J134").Value = 1007
$ws.Range("K134").Value = 3028.5
$ws.Range("L134").Value = 3021
$ws.Range("M134").Value = -493.5
$ws.Range("N134").Value = -8091

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7442.2383
$ws.Range("I31").Value = 3066.5
$ws.Range("J31").Value = 9192.532999999999
$ws.Range("K31").Value = 3066.5
$ws.Range("L31").Value = 9192.532999999999
$ws.Range("M31").Value = -2771.5
$ws.Range("N31").Value = -9782.532999999999

$ws.Range("H34").Value = 7442.2383
$ws.Range("I34").Value = 3066.5
$ws.Range("J34").Value = 9192.532999999999
$ws.Range("K34").Value = 3066.5
$ws.Range("L34").Value = 9192.532999999999
$ws.Range("M34").Value = -2864.5
$ws.Range("N34").Value = -9596.532999999999

$ws.Range("H58").Value = 7068.4287
$ws.Range("I58").Value = 4000
$ws.Range("J58").Value = 7579.8335
$ws.Range("K58").Value = 4000
$ws.Range("L58").Value = 7579.8335
$ws.Range("M58").Value = -3797
$ws.Range("N58").Value = -7985.8335

$ws.Range("H62").Value = 7714.143
$ws.Range("J62").Value = 9250
$ws.Range("L62").Value = 9250
$ws.Range("N62").Value = -10498

$ws.Range("H65").Value = 7714.143
$ws.Range("J65").Value = 9250
$ws.Range("L65").Value = 46250
$ws.Range("N65").Value = -52490

$ws.Range("H123").Value = 48983.332
$ws.Range("J123").Value = 48983.332
$ws.Range("L123").Value = 48983.332
$ws.Range("N123").Value = -58783.332

$ws.Range("H134").Value = 1132.4166
$ws.Range("I134").Value = 1199.4445
$ws.Range("J134").Value = 931.3333
$ws.Range("K134").Value = 3598.3335
$ws.Range("L134").Value = 2793.9999
$ws.Range("M134").Value = -1063.3335
$ws.Range("N134").Value = -7863.9999

$ws.Range("H136").Value = 7068.4287
$ws.Range("I136").Value = 4000
$ws.Range("J136").Value = 7579.8335
$ws.Range("K136").Value = 12000
$ws.Range("L136").Value = 22739.5005
$ws.Range("M136").Value = -9450
$ws.Range("N136").Value = -27839.5005

$ws.Range("H140").Value = 75780
$ws.Range("J140").Value = 75780
$ws.Range("L140").Value = 75780
$ws.Range("N140").Value = -86140

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 5360
$ws.Range("J137").Value = 5231.6665
$ws.Range("L137").Value = 15694.9995
$ws.Range("N137").Value = -25894.9995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2990.125
$ws.Range("I102").Value = 1902.2
$ws.Range("J102").Value = 4803.3335
$ws.Range("K102").Value = 1902.2
$ws.Range("L102").Value = 4803.3335
$ws.Range("M102").Value = -280.2
$ws.Range("N102").Value = -8047.3335

$ws.Range("H122").Value = 4350.5
$ws.Range("I122").Value = 3563.125
$ws.Range("K122").Value = 10689.375
$ws.Range("M122").Value = -8239.375

$ws.Range("H140").Value = 101391.86
$ws.Range("J140").Value = 77867.39999999999
$ws.Range("L140").Value = 77867.39999999999
$ws.Range("N140").Value = -88227.39999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2342.5652
$ws.Range("I40").Value = 2269.1
$ws.Range("J40").Value = 2832.3333
$ws.Range("K40").Value = 2269.1
$ws.Range("L40").Value = 2832.3333
$ws.Range("M40").Value = -2133.1
$ws.Range("N40").Value = -3104.3333

$ws.Range("H68").Value = 5900
$ws.Range("I68").Value = 2500
$ws.Range("K68").Value = 2500
$ws.Range("M68").Value = -1751

$ws.Range("H71").Value = 5900
$ws.Range("I71").Value = 2500
$ws.Range("K71").Value = 12500
$ws.Range("M71").Value = -8756

$ws.Range("H132").Value = 2986.5
$ws.Range("I132").Value = 2698.8572
$ws.Range("K132").Value = 8096.571599999999
$ws.Range("M132").Value = -5566.571599999999

$ws.Range("H136").Value = 4548.8184
$ws.Range("I136").Value = 4421.5713
$ws.Range("J136").Value = 4771.5
$ws.Range("K136").Value = 13264.7139
$ws.Range("L136").Value = 14314.5
$ws.Range("M136").Value = -10714.7139
$ws.Range("N136").Value = -19414.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2520
$ws.Range("J136").Value = 4366.5
$ws.Range("L136").Value = 13099.5
$ws.Range("N136").Value = -18199.5
